$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that were previously blank and now get the text "NULL" (red font,
# reusing the existing shared string / cell style already present in the
# workbook).
$cells = @("D5", "G10", "B16", "C16", "E22", "D30", "E32", "C33", "B48", "C48", "B55", "G57", "C58")

foreach ($cellRef in $cells) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "NULL"
    $rng.Font.Color = 255
}

# Restore the last active selection recorded in the saved sheet view.
$ws.Range("G10").Select() | Out-Null
